# Update extrapolation calibration values to remove noisy sub-$5 price rows
# (commit: "Removing less than USD 5 price from extrapolation calibration
# because it is just a noise")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7
$ws.Range("D7").Value = 121208.39890571
$ws.Range("E7").Value = -0.02675736004434551
$ws.Range("F7").Value = 0.2427200386657377
$ws.Range("G7").Value = -0.7772775658794857
$ws.Range("H7").Value = 5.878484890207581

# Row 8
$ws.Range("D8").Value = 121473.2495089783
$ws.Range("E8").Value = -0.04344350487129781
$ws.Range("F8").Value = 0.2090703374888578
$ws.Range("G8").Value = -0.8501621764412184
$ws.Range("H8").Value = 6.658461103396522

# Row 9
$ws.Range("D9").Value = 123565.7903413193
$ws.Range("E9").Value = -0.07364479176811603
$ws.Range("F9").Value = 0.3595636893797028
$ws.Range("G9").Value = -2.022994663812912
$ws.Range("H9").Value = 12.75122567204191

# Row 10
$ws.Range("D10").Value = 124754.5125494682
$ws.Range("E10").Value = -0.1118687796455338
$ws.Range("F10").Value = 0.4390188915677648
$ws.Range("G10").Value = -1.926370987563858
$ws.Range("H10").Value = 9.911098807693437

# Row 13
$ws.Range("D13").Value = 116976.6969166577
$ws.Range("E13").Value = 0.004491494622911664
$ws.Range("F13").Value = 0.1009526540357622
$ws.Range("G13").Value = -0.5829877988252592
$ws.Range("H13").Value = 9.521702830545536

# Row 17
$ws.Range("D17").Value = 116917.6199590839
$ws.Range("E17").Value = 0.006169810849410319
$ws.Range("F17").Value = 0.08645647194491804
$ws.Range("G17").Value = -0.7421589743455126
$ws.Range("H17").Value = 7.279101729407761

# Row 18
$ws.Range("D18").Value = 116962.3750568939
$ws.Range("E18").Value = -0.0001193861189121011
$ws.Range("F18").Value = 0.1063699877475345
$ws.Range("G18").Value = -0.5289535724714836
$ws.Range("H18").Value = 6.748736341346964

# Row 20
$ws.Range("D20").Value = 117776.3619845627
$ws.Range("E20").Value = 0.007624718392454
$ws.Range("F20").Value = 0.1342423211913152
$ws.Range("G20").Value = -0.2301401050947217
$ws.Range("H20").Value = 5.815856764522648

$wb.Save()
